# Updated data to reflect new requirement separation
# Insert three new columns (Corequisites, Concurrent, Recommended) between
# the existing "Prerequisites" (C) and "Terms Typically Offered" (old D, now G)
# columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank columns at D:F. This shifts the existing "Terms
# Typically Offered" column from D to G, carrying its data along.
$ws.Columns("D:F").Insert()

# New header row values for the inserted columns (G1 already holds the
# shifted-over "Terms Typically Offered" header).
$ws.Range("D1").Value2 = "Corequisites"
$ws.Range("E1").Value2 = "Concurrent"
$ws.Range("F1").Value2 = "Recommended"

# Default all new Corequisites/Concurrent/Recommended cells (rows 2-36) to "NA".
$ws.Range("D2:F36").Value2 = "NA"

# Non-breaking space used inside some course-number references (e.g. "TH 210")
$nbsp = [char]0x00A0

# Rows 18, 24 and 27 share the same "Junior standing; ... Recommended: TH 210
# for all other majors." prerequisite text. The trailing "Recommended: ..."
# clause is being split out into the new "Recommended" column, so strip it
# from the Prerequisites text and place it into column F instead. The
# "Terms Typically Offered" values for these three rows also pick up a
# trailing space in the new layout.
$newPrereq = "Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and TH${nbsp}210, TH${nbsp}227, or TH${nbsp}228 for Theatre majors or completion of GE Area C3 for all other majors."
$recommendedText = "TH${nbsp}210 for all other majors."

$specialRows = @(
    @{ Row = 18; Terms = "W " },
    @{ Row = 24; Terms = "TBD " },
    @{ Row = 27; Terms = "SU " }
)

foreach ($item in $specialRows) {
    $r = $item.Row
    $ws.Range("C$r").Value2 = $newPrereq
    $ws.Range("F$r").Value2 = $recommendedText
    $ws.Range("G$r").Value2 = $item.Terms
}
